$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Internal Assignment" column (column O) contents in the
# sample-type export sheet: header in O4 and the per-row values in O5:O7.
$ws.Range("O4:O7").ClearContents()
